$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows corresponding to Neutrophils and Resolving-Mac target clusters
# (for both ECs and FAPs sending clusters). Delete bottom-up to keep row numbers valid.
$ws.Rows.Item(13).Delete()  # FAPs -> Resolving-Mac
$ws.Rows.Item(12).Delete()  # FAPs -> Neutrophils
$ws.Rows.Item(7).Delete()   # ECs -> Resolving-Mac
$ws.Rows.Item(6).Delete()   # ECs -> Neutrophils

# Update the recomputed TPM-based statistics for the remaining rows
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.5
$ws.Range("G2").Value2 = 0.0292345
$ws.Range("H2").Value2 = 0.058469
$ws.Range("I2").Value2 = 0.4428765120700495
$ws.Range("J2").Value2 = 0.346386487911515
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 30.801072
$ws.Range("N2").Value2 = 61.602144
$ws.Range("O2").Value2 = 0.5373480691764108
$ws.Range("P2").Value2 = 0.485871843331092
$ws.Range("Q2").Value2 = 0.900453939384
$ws.Range("R2").Value2 = 3.601815757536
$ws.Range("S2").Value2 = 0.2379788386444245
$ws.Range("T2").Value2 = 0.1682994413865508
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.5
$ws.Range("G3").Value2 = 0.0292345
$ws.Range("H3").Value2 = 0.058469
$ws.Range("I3").Value2 = 0.4428765120700495
$ws.Range("J3").Value2 = 0.346386487911515
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 12.015213
$ws.Range("N3").Value2 = 36.045639
$ws.Range("O3").Value2 = 0.2096145064786482
$ws.Range("P3").Value2 = 0.2843011610923331
$ws.Range("Q3").Value2 = 0.3512587444485
$ws.Range("R3").Value2 = 2.107552466691
$ws.Range("S3").Value2 = 0.09283334150854852
$ws.Range("T3").Value2 = 0.09847808069993912
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.5
$ws.Range("G4").Value2 = 0.0292345
$ws.Range("H4").Value2 = 0.058469
$ws.Range("I4").Value2 = 0.4428765120700495
$ws.Range("J4").Value2 = 0.346386487911515
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.1305583333333333
$ws.Range("N4").Value2 = 0.391675
$ws.Range("O4").Value2 = 0.002277689176907768
$ws.Range("P4").Value2 = 0.003089240761436898
$ws.Range("Q4").Value2 = 0.003816807595833333
$ws.Range("R4").Value2 = 0.022900845575
$ws.Range("S4").Value2 = 0.001008735038248614
$ws.Range("T4").Value2 = 0.001070071257667222
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.5
$ws.Range("G5").Value2 = 0.0292345
$ws.Range("H5").Value2 = 0.058469
$ws.Range("I5").Value2 = 0.4428765120700495
$ws.Range("J5").Value2 = 0.346386487911515
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 14.373679
$ws.Range("N5").Value2 = 28.747358
$ws.Range("O5").Value2 = 0.2507597351680332
$ws.Range("P5").Value2 = 0.2267377548151379
$ws.Range("Q5").Value2 = 0.4202073187255
$ws.Range("R5").Value2 = 1.680829274902
$ws.Range("S5").Value2 = 0.1110555968788279
$ws.Range("T5").Value2 = 0.0785388945673578
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.036776
$ws.Range("H6").Value2 = 0.110328
$ws.Range("I6").Value2 = 0.5571234879299505
$ws.Range("J6").Value2 = 0.6536135120884849
$ws.Range("K6").Value2 = 2
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 30.801072
$ws.Range("N6").Value2 = 61.602144
$ws.Range("O6").Value2 = 0.5373480691764108
$ws.Range("P6").Value2 = 0.485871843331092
$ws.Range("Q6").Value2 = 1.132740223872
$ws.Range("R6").Value2 = 6.796441343231999
$ws.Range("S6").Value2 = 0.2993692305319863
$ws.Range("T6").Value2 = 0.3175724019445412
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.036776
$ws.Range("H7").Value2 = 0.110328
$ws.Range("I7").Value2 = 0.5571234879299505
$ws.Range("J7").Value2 = 0.6536135120884849
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 12.015213
$ws.Range("N7").Value2 = 36.045639
$ws.Range("O7").Value2 = 0.2096145064786482
$ws.Range("P7").Value2 = 0.2843011610923331
$ws.Range("Q7").Value2 = 0.441871473288
$ws.Range("R7").Value2 = 3.976843259592
$ws.Range("S7").Value2 = 0.1167811649700997
$ws.Range("T7").Value2 = 0.185823080392394
$ws.Range("E8").Value2 = 1
$ws.Range("F8").Value2 = 0.3333333333333333
$ws.Range("G8").Value2 = 0.036776
$ws.Range("H8").Value2 = 0.110328
$ws.Range("I8").Value2 = 0.5571234879299505
$ws.Range("J8").Value2 = 0.6536135120884849
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.1305583333333333
$ws.Range("N8").Value2 = 0.391675
$ws.Range("O8").Value2 = 0.002277689176907768
$ws.Range("P8").Value2 = 0.003089240761436898
$ws.Range("Q8").Value2 = 0.004801413266666666
$ws.Range("R8").Value2 = 0.0432127194
$ws.Range("S8").Value2 = 0.001268954138659154
$ws.Range("T8").Value2 = 0.002019169503769676
$ws.Range("E9").Value2 = 1
$ws.Range("F9").Value2 = 0.3333333333333333
$ws.Range("G9").Value2 = 0.036776
$ws.Range("H9").Value2 = 0.110328
$ws.Range("I9").Value2 = 0.5571234879299505
$ws.Range("J9").Value2 = 0.6536135120884849
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 14.373679
$ws.Range("N9").Value2 = 28.747358
$ws.Range("O9").Value2 = 0.2507597351680332
$ws.Range("P9").Value2 = 0.2267377548151379
$ws.Range("Q9").Value2 = 0.528606418904
$ws.Range("R9").Value2 = 3.171638513424
$ws.Range("S9").Value2 = 0.1397041382892053
$ws.Range("T9").Value2 = 0.14819886024778
